# "finish 1st draft Ch 7" - text fixes + strike-through annotations.

$d = $word.ActiveDocument

# --- Ch 6 (Managing AD) intro: British -> US spelling ("organisations" -> "organizations") ---
$d.Content.Find.Execute(
    "organisations these days", $true, $false, $false, $false, $false,
    $true, 1, $false, "organizations these days", 2) | Out-Null

# --- Ch 7 (Enterprise security) intro paragraph fixes ---
# 1) "how to establish implement" -> "how to establish and implement"
$d.Content.Find.Execute(
    "how to establish implement fine-grained", $true, $false, $false, $false, $false,
    $true, 1, $false, "how to establish and implement fine-grained", 2) | Out-Null

# 2) "how to install/use AD based certificate" -> "how you install/use AD-based certificate"
$d.Content.Find.Execute(
    "fine-grained security delegation, how to install/use AD based certificate",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "fine-grained security delegation, how you install/use AD-based certificate", 2) | Out-Null

# --- "Privileged users" bullet: append a struck-through annotation ---
$r = $d.Content
$r.Find.Execute("Privileged users", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter("  (done in earlier chapter)")
$r.Font.StrikeThrough = 1

# --- "Installing WSL and WSL 2" bullet: append a struck-through annotation ---
$r2 = $d.Content
$r2.Find.Execute("Installing WSL and WSL 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$r2.InsertAfter(" (not needed).")
$r2.Font.StrikeThrough = 1

# --- Ch 13 (Managing Hyper-V) intro: British -> US spelling ("virtualisation" -> "virtualization") ---
$d.Content.Find.Execute(
    "native virtualisation offering", $true, $false, $false, $false, $false,
    $true, 1, $false, "native virtualization offering", 2) | Out-Null
